$d = $word.ActiveDocument

# Remove the run containing "Try sth" (the whole run/text), leaving the
# paragraph mark and the bookmark that follows it intact.
$d.Content.Find.Execute("Try sth", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
